$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titleText = "Play Candy Bars Free - Review of IGT's Slot Game"
$metaLabel = "Meta description"
$metaRest  = ": Read our review of Candy Bars by IGT. Play this colorful and classic slot game for free. Learn how to win Blackout Wins and Progressive Jackpots."

$titlePara = $d.Paragraphs(1)

# Create a brand new (empty) paragraph straight after the title paragraph.
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# Fill it with the full sentence, then go back and bold just the label part.
$metaPara.Range.Text = $metaLabel + $metaRest
$labelRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $metaLabel.Length)
$labelRange.Bold = 1

# Everything from here on should only look *after* the block we just
# inserted, so the searches below can't accidentally match the text we just
# typed in (it happens to reuse the same sentence).
$tailStart = $metaPara.Range.End

# ---------------------------------------------------------------------------
# 2) Drop the duplicated bold title paragraph near the end of the document.
# ---------------------------------------------------------------------------
$dupScope = $d.Range($tailStart, $d.Content.End)
$dupScope.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($dupScope.Find.Found) {
    $dupPara = $dupScope.Paragraphs(1)
    $dupPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the old italic meta-description text with the new image prompt.
# ---------------------------------------------------------------------------
$oldPrompt = "Read our review of Candy Bars by IGT. Play this colorful and classic slot game for free. Learn how to win Blackout Wins and Progressive Jackpots."
$newPrompt = "Create a vibrant feature image for Candy Bars that features a happy Maya warrior wearing glasses, in a cartoon style. The background should be bright and colorful, with a candy-themed design such as candy canes, gumdrops, and lollipops. The Maya warrior should be holding a big lollipop and have a big smile on their face, with candy symbols surrounding them such as gumballs and chocolate bars. The image should showcase the fun and playful nature of the game while incorporating its candy theme and the idea of winning big."

$promptScope = $d.Range($tailStart, $d.Content.End)
$promptScope.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 1) | Out-Null
